$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2-7 and write new rows 8-10 to reflect the added "ECs" sending-cluster group
# Row 2: ECs -> ECs
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Inhba"
$ws.Cells.Item(2, 3).Value = "Acvr1b"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 3.675031333333333
$ws.Cells.Item(2, 8).Value = 11.025094
$ws.Cells.Item(2, 9).Value = 0.2032371147293133
$ws.Cells.Item(2, 10).Value = 0.2032371147293133
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 4.021200333333334
$ws.Cells.Item(2, 14).Value = 12.063601
$ws.Cells.Item(2, 15).Value = 0.389801966361343
$ws.Cells.Item(2, 16).Value = 0.389801966361343
$ws.Cells.Item(2, 17).Value = 14.77803722261045
$ws.Cells.Item(2, 18).Value = 133.002335003494
$ws.Cells.Item(2, 19).Value = 0.07922222695909219
$ws.Cells.Item(2, 20).Value = 0.07922222695909219

# Row 3: ECs -> FAPs
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Inhba"
$ws.Cells.Item(3, 3).Value = "Acvr1b"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 3.675031333333333
$ws.Cells.Item(3, 8).Value = 11.025094
$ws.Cells.Item(3, 9).Value = 0.2032371147293133
$ws.Cells.Item(3, 10).Value = 0.2032371147293133
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 3.562995333333333
$ws.Cells.Item(3, 14).Value = 10.688986
$ws.Cells.Item(3, 15).Value = 0.3453850770768087
$ws.Cells.Item(3, 16).Value = 0.3453850770768087
$ws.Cells.Item(3, 17).Value = 13.09411949052044
$ws.Cells.Item(3, 18).Value = 117.847075414684
$ws.Cells.Item(3, 19).Value = 0.07019506653565208
$ws.Cells.Item(3, 20).Value = 0.07019506653565208

# Row 4: ECs -> sCs
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Inhba"
$ws.Cells.Item(4, 3).Value = "Acvr1b"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 7).Value = 3.675031333333333
$ws.Cells.Item(4, 8).Value = 11.025094
$ws.Cells.Item(4, 9).Value = 0.2032371147293133
$ws.Cells.Item(4, 10).Value = 0.2032371147293133
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 2.731812666666666
$ws.Cells.Item(4, 14).Value = 8.195438
$ws.Cells.Item(4, 15).Value = 0.2648129565618484
$ws.Cells.Item(4, 16).Value = 0.2648129565618484
$ws.Cells.Item(4, 17).Value = 10.03949714679689
$ws.Cells.Item(4, 18).Value = 90.35547432117198
$ws.Cells.Item(4, 19).Value = 0.05381982123456905
$ws.Cells.Item(4, 20).Value = 0.05381982123456905

# Row 5: FAPs -> ECs
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Inhba"
$ws.Cells.Item(5, 3).Value = "Acvr1b"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 10.108494
$ws.Cells.Item(5, 8).Value = 30.325482
$ws.Cells.Item(5, 9).Value = 0.5590213983169419
$ws.Cells.Item(5, 10).Value = 0.5590213983169419
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 4.021200333333334
$ws.Cells.Item(5, 14).Value = 12.063601
$ws.Cells.Item(5, 15).Value = 0.389801966361343
$ws.Cells.Item(5, 16).Value = 0.389801966361343
$ws.Cells.Item(5, 17).Value = 40.648279442298
$ws.Cells.Item(5, 18).Value = 365.834514980682
$ws.Cells.Item(5, 19).Value = 0.2179076403020115
$ws.Cells.Item(5, 20).Value = 0.2179076403020115

# Row 6: FAPs -> FAPs
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Inhba"
$ws.Cells.Item(6, 3).Value = "Acvr1b"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 10.108494
$ws.Cells.Item(6, 8).Value = 30.325482
$ws.Cells.Item(6, 9).Value = 0.5590213983169419
$ws.Cells.Item(6, 10).Value = 0.5590213983169419
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 3.562995333333333
$ws.Cells.Item(6, 14).Value = 10.688986
$ws.Cells.Item(6, 15).Value = 0.3453850770768087
$ws.Cells.Item(6, 16).Value = 0.3453850770768087
$ws.Cells.Item(6, 17).Value = 36.016516949028
$ws.Cells.Item(6, 18).Value = 324.148652541252
$ws.Cells.Item(6, 19).Value = 0.1930776487452823
$ws.Cells.Item(6, 20).Value = 0.1930776487452823

# Row 7: FAPs -> sCs
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Inhba"
$ws.Cells.Item(7, 3).Value = "Acvr1b"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 10.108494
$ws.Cells.Item(7, 8).Value = 30.325482
$ws.Cells.Item(7, 9).Value = 0.5590213983169419
$ws.Cells.Item(7, 10).Value = 0.5590213983169419
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 2.731812666666666
$ws.Cells.Item(7, 14).Value = 8.195438
$ws.Cells.Item(7, 15).Value = 0.2648129565618484
$ws.Cells.Item(7, 16).Value = 0.2648129565618484
$ws.Cells.Item(7, 17).Value = 27.614511950124
$ws.Cells.Item(7, 18).Value = 248.530607551116
$ws.Cells.Item(7, 19).Value = 0.1480361092696481
$ws.Cells.Item(7, 20).Value = 0.1480361092696481

# Row 8: sCs -> ECs
$ws.Cells.Item(8, 1).Value = "sCs"
$ws.Cells.Item(8, 2).Value = "Inhba"
$ws.Cells.Item(8, 3).Value = "Acvr1b"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 4.298956
$ws.Cells.Item(8, 8).Value = 12.896868
$ws.Cells.Item(8, 9).Value = 0.2377414869537448
$ws.Cells.Item(8, 10).Value = 0.2377414869537448
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 4.021200333333334
$ws.Cells.Item(8, 14).Value = 12.063601
$ws.Cells.Item(8, 15).Value = 0.389801966361343
$ws.Cells.Item(8, 16).Value = 0.389801966361343
$ws.Cells.Item(8, 17).Value = 17.28696330018534
$ws.Cells.Item(8, 18).Value = 155.582669701668
$ws.Cells.Item(8, 19).Value = 0.0926720991002393
$ws.Cells.Item(8, 20).Value = 0.09267209910023928

# Row 9: sCs -> FAPs
$ws.Cells.Item(9, 1).Value = "sCs"
$ws.Cells.Item(9, 2).Value = "Inhba"
$ws.Cells.Item(9, 3).Value = "Acvr1b"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 4.298956
$ws.Cells.Item(9, 8).Value = 12.896868
$ws.Cells.Item(9, 9).Value = 0.2377414869537448
$ws.Cells.Item(9, 10).Value = 0.2377414869537448
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 3.562995333333333
$ws.Cells.Item(9, 14).Value = 10.688986
$ws.Cells.Item(9, 15).Value = 0.3453850770768087
$ws.Cells.Item(9, 16).Value = 0.3453850770768087
$ws.Cells.Item(9, 17).Value = 15.31716016620534
$ws.Cells.Item(9, 18).Value = 137.854441495848
$ws.Cells.Item(9, 19).Value = 0.08211236179587424
$ws.Cells.Item(9, 20).Value = 0.08211236179587424

# Row 10: sCs -> sCs
$ws.Cells.Item(10, 1).Value = "sCs"
$ws.Cells.Item(10, 2).Value = "Inhba"
$ws.Cells.Item(10, 3).Value = "Acvr1b"
$ws.Cells.Item(10, 4).Value = "sCs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 4.298956
$ws.Cells.Item(10, 8).Value = 12.896868
$ws.Cells.Item(10, 9).Value = 0.2377414869537448
$ws.Cells.Item(10, 10).Value = 0.2377414869537448
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 2.731812666666666
$ws.Cells.Item(10, 14).Value = 8.195438
$ws.Cells.Item(10, 15).Value = 0.2648129565618484
$ws.Cells.Item(10, 16).Value = 0.2648129565618484
$ws.Cells.Item(10, 17).Value = 11.74394245424267
$ws.Cells.Item(10, 18).Value = 105.695482088184
$ws.Cells.Item(10, 19).Value = 0.06295702605763127
$ws.Cells.Item(10, 20).Value = 0.06295702605763127
